# Apply edit: rename sheet, insert 4 header rows with master-package info,
# and resize columns C/D now that the tolerance columns moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet/tab from "Page 1 - Table 1" to "Sheet1"
$ws.Name = "Sheet1"

# 2) Insert 4 new rows above the current row 1, pushing everything else down.
$insertRange = $ws.Range("A1:J4")
$insertRange.EntireRow.Insert()

# 3) Populate the new header rows with the master package info.
$ws.Range("A1").Value = "MASTER PACKAGE"
$ws.Range("A2").Value = "WesternGlove Centric8 PROD"
$ws.Range("B2").Value = "M12225BVS563:KONRAD"
$ws.Range("C2").Value = "DUP REVIEW"
$ws.Range("D2").Value = "Revised 1/9/25, 2:47 PM"
$ws.Range("A3").Value = "Evaluation"

# 4) Re-apply the sheet-wide style (thin border, left/top aligned, wrap text)
#    to the new cells so they match the rest of the table's look.
$styleRange = $ws.Range("A1:D3")
$styleRange.Borders.LineStyle = 1
$styleRange.Borders.Weight = 2
$styleRange.HorizontalAlignment = -4131
$styleRange.VerticalAlignment = -4160
$styleRange.WrapText = $true

# Also keep the now-blank cells (B1,C1,D1,B3,C3,D3) consistently styled.
$ws.Range("B1:D1").Borders.LineStyle = 1
$ws.Range("B3:D3").Borders.LineStyle = 1

# 5) Resize columns C and D now that they no longer share one narrow width.
$ws.Columns("C").ColumnWidth = 12
$ws.Columns("D").ColumnWidth = 22
